# Commit: "changing document, table attributes to lowerCamelCase"
#
# The workbook uses a lightweight text-based markup (ObjTables) embedded in
# inline-string cells of each worksheet to describe table metadata
# (version / type / id). This change renames the markup's attribute
# keywords from UpperCamelCase to lowerCamelCase:
#   ObjTablesVersion -> objTablesVersion
#   Type             -> type
#   Id               -> id

$wb = $excel.ActiveWorkbook

$ws_mainroot = $wb.Worksheets.Item("!!Main root")
$ws_mainroot.Range("A1").Value = "!!!ObjTables objTablesVersion='0.0.8'"
$ws_mainroot.Range("A2").Value = "!!ObjTables type='Data' id='MainRoot'"

$ws_nodes = $wb.Worksheets.Item("!!Nodes")
$ws_nodes.Range("A1").Value = "!!ObjTables type='Data' id='Node'"

$ws_nodefriends = $wb.Worksheets.Item("!!Node friends")
$ws_nodefriends.Range("A1").Value = "!!ObjTables type='Data' id='NodeFriend'"

$ws_leaves = $wb.Worksheets.Item("!!Leaves")
$ws_leaves.Range("A1").Value = "!!ObjTables type='Data' id='Leaf'"

$ws_rows = $wb.Worksheets.Item("!!One to many rows")
$ws_rows.Range("A1").Value = "!!ObjTables type='Data' id='OneToManyRow'"
